$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 12:22"

# Rumania overtakes Peru in total cases (5202 -> 5467), swapping their
# ranking rows (33/34). Peru's own figures are unchanged, just shifted
# down one row.
$ws.Range("A33").Value = "Rumania"
$ws.Range("B33").Value = 5467
$ws.Range("C33").Value = 265
$ws.Range("D33").Value = 729
$ws.Range("E33").Value = 4481
$ws.Range("F33").Value = 183
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 257

$ws.Range("A34").Value = "Peru"
$ws.Range("B34").Value = 5256
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 1438
$ws.Range("E34").Value = 3680
$ws.Range("F34").Value = 124
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 138

# Catar overtakes Tailandia in total cases (2376 -> 2512), swapping their
# ranking rows (47/48). Tailandia's own figures are unchanged, just
# shifted down one row.
$ws.Range("A47").Value = "Catar"
$ws.Range("B47").Value = 2512
$ws.Range("C47").Value = 136
$ws.Range("D47").Value = 227
$ws.Range("E47").Value = 2279
$ws.Range("F47").Value = 37
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 6

$ws.Range("A48").Value = "Tailandia"
$ws.Range("B48").Value = 2473
$ws.Range("C48").Value = 50
$ws.Range("D48").Value = 1013
$ws.Range("E48").Value = 1427
$ws.Range("F48").Value = 61
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 33

# Eslovenia figures refreshed (no rank change)
$ws.Range("B67").Value = 1160
$ws.Range("C67").Value = 36
$ws.Range("D67").Value = 137
$ws.Range("E67").Value = 978
$ws.Range("F67").Value = 36
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 45

# Hong Kong figures refreshed (no rank change)
$ws.Range("B70").Value = 990
$ws.Range("C70").Value = 16
$ws.Range("D70").Value = 309
$ws.Range("E70").Value = 677
$ws.Range("F70").Value = 15
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 4

# Etiopia: "Muertes hoy" (F) corrected 2 -> 1
$ws.Range("F140").Value = 1
